$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Strip the "gi|<id>|ref|" prefix and the trailing "|" from the accession
# strings in column B, leaving just "<accession>-<position>".
# Values are written in the same order the original strings were first
# introduced into the shared-string table so the table layout lines up
# with the target workbook.
$ws.Range("B2").Value = "NC_017251.1-213522"
$ws.Range("B13").Value = "NC_017251.1-1146371"
$ws.Range("B18").Value = "NC_017251.1-549890"
$ws.Range("B20").Value = "NC_017251.1-1184227"
$ws.Range("B21").Value = "NC_017251.1-1245156"
$ws.Range("B25").Value = "NC_017251.1-994359"
$ws.Range("B27").Value = "NC_017251.1-1088654"
$ws.Range("B3").Value = "NC_017250.1-1173757"
$ws.Range("B4").Value = "NC_017250.1-241603"
$ws.Range("B5").Value = "NC_017250.1-1013589"
$ws.Range("B6").Value = "NC_017250.1-1025446"
$ws.Range("B7").Value = "NC_017250.1-1020241"
$ws.Range("B8").Value = "NC_017250.1-1175324"
$ws.Range("B9").Value = "NC_017250.1-291907"
$ws.Range("B10").Value = "NC_017250.1-1197913"
$ws.Range("B11").Value = "NC_017250.1-1016576"
$ws.Range("B12").Value = "NC_017250.1-264518"
$ws.Range("B14").Value = "NC_017250.1-1014292"
$ws.Range("B15").Value = "NC_017250.1-1048661"
$ws.Range("B16").Value = "NC_017250.1-1072841"
$ws.Range("B17").Value = "NC_017250.1-1099151"
$ws.Range("B19").Value = "NC_017250.1-1167451"
$ws.Range("B22").Value = "NC_017250.1-1011109"
$ws.Range("B23").Value = "NC_017250.1-1035261"
$ws.Range("B24").Value = "NC_017250.1-1147663"
$ws.Range("B26").Value = "NC_017250.1-1195126"

# The accession/position strings are free-form text (some start with
# digits), so format column B as Text to keep them from being
# reinterpreted.
$ws.Range("B2:B27").NumberFormat = "@"

# Move the active selection to match the saved workbook state.
$ws.Range("K25").Select()
